# Auto-generated edit script: apply recalculated output values after
# setting "Unmet Demand Penalty" (Summary!B3) to 0.01.
# The workbook stores only static cached values (no formulas), so each
# affected cell is written directly with its new value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = 131579.9112553819
$ws.Range("B8").Value = 24253065.61257719
$ws.Range("B10").Value = 2491228.976683192

$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G5").Value = 414.4337959369544
$ws.Range("H5").Value = 330.5757541782243
$ws.Range("I5").Value = 176.9760193775952
$ws.Range("J5").Value = 107.2955742555736
$ws.Range("K5").Value = 109.5572237694796
$ws.Range("L5").Value = 98.64091687123928
$ws.Range("M5").Value = 77.7676953375541
$ws.Range("N5").Value = 74.3656454478664
$ws.Range("O5").Value = 83.69133109099639
$ws.Range("P5").Value = 106.2781106359148
$ws.Range("Q5").Value = 128.4697750236904
$ws.Range("R5").Value = 161.0018864037399
$ws.Range("S5").Value = 189.2190633734531
$ws.Range("T5").Value = 219.2920578056454
$ws.Range("U5").Value = 251.276137581582
$ws.Range("G6").Value = 136.8785924310737
$ws.Range("H6").Value = 107.7452501129632
$ws.Range("I6").Value = 83.51510303826707
$ws.Range("J6").Value = 82.91243519753434
$ws.Range("K6").Value = 62.76629045205057
$ws.Range("L6").Value = 37.60657774285653
$ws.Range("M6").Value = 24.33271034503603
$ws.Range("N6").Value = 10.42253800004659
$ws.Range("O6").Value = 31.97882363640291
$ws.Range("P6").Value = 45.19417501179163
$ws.Range("Q6").Value = 80.63453985745144
$ws.Range("R6").Value = 116.8133877002326
$ws.Range("S6").Value = 163.0473981187501
$ws.Range("T6").Value = 198.2907557613397
$ws.Range("U6").Value = 225.9107949275447
$ws.Range("G7").Value = 167.6012020808691
$ws.Range("H7").Value = 158.7616981666879
$ws.Range("I7").Value = 143.7288091611985
$ws.Range("J7").Value = 99.43816791380756
$ws.Range("K7").Value = 83.72822537421013
$ws.Range("L7").Value = 76.93542539304551
$ws.Range("M7").Value = 77.82642397052864
$ws.Range("N7").Value = 68.03899070462725
$ws.Range("O7").Value = 83.36329197944329
$ws.Range("P7").Value = 90.5862140395771
$ws.Range("Q7").Value = 118.8664697760067
$ws.Range("R7").Value = 159.7675876048201
$ws.Range("S7").Value = 217.2238431175235
$ws.Range("T7").Value = 226.2801774240348
$ws.Range("U7").Value = 286.2977687777133
$ws.Range("G8").Value = 414.2543098065221
$ws.Range("H8").Value = 328.737591844935
$ws.Range("I8").Value = 170.0563803341062
$ws.Range("J8").Value = 92.06191329280011
$ws.Range("K8").Value = 86.72591490550931
$ws.Range("L8").Value = 70.31665934305323
$ws.Range("M8").Value = 46.25150133729363
$ws.Range("N8").Value = 42.33948647951536
$ws.Range("O8").Value = 53.44993733213093
$ws.Range("P8").Value = 80.46778072209531
$ws.Range("Q8").Value = 109.0872921559756
$ws.Range("R8").Value = 149.7272407629755
$ws.Range("S8").Value = 185.1290231762283
$ws.Range("T8").Value = 218.5063572696783
$ws.Range("U8").Value = 251.2617786911475
$ws.Range("G9").Value = 136.7825588849655
$ws.Range("H9").Value = 106.817768233445
$ws.Range("I9").Value = 80.20868489375378
$ws.Range("J9").Value = 73.83937108913167
$ws.Range("K9").Value = 47.25897875440141
$ws.Range("L9").Value = 16.75508344424118
$ws.Range("M9").Value = 0
$ws.Range("O9").Value = 9.130000058643361
$ws.Range("P9").Value = 26.85597970277001
$ws.Range("Q9").Value = 68.37594193669796
$ws.Range("R9").Value = 110.8508838462542
$ws.Range("S9").Value = 161.2636171197674
$ws.Range("T9").Value = 197.9036731785617
$ws.Range("U9").Value = 225.9044769310903
$ws.Range("G10").Value = 167.5206907817879
$ws.Range("H10").Value = 158.0458795257661
$ws.Range("I10").Value = 141.3076148215571
$ws.Range("J10").Value = 93.74601906876785
$ws.Range("K10").Value = 74.37427626277801
$ws.Range("L10").Value = 64.96559098237566
$ws.Range("M10").Value = 65.20591187910269
$ws.Range("N10").Value = 55.71856618250388
$ws.Range("O10").Value = 71.98338581476828
$ws.Range("P10").Value = 80.84873837615825
$ws.Range("Q10").Value = 112.1247463593087
$ws.Range("R10").Value = 156.1475068297698
$ws.Range("S10").Value = 215.8207507508087
$ws.Range("T10").Value = 225.9361746006879
$ws.Range("U10").Value = 286.2933772523089
$ws.Range("L11").Value = 28.83625887535973
$ws.Range("M11").Value = 0.09656339947139259
$ws.Range("L13").Value = 47.43596908875769
$ws.Range("N13").Value = 37.67551094615368
$ws.Range("L17").Value = 28.8362588753597
$ws.Range("M17").Value = 0.09656339947136416
$ws.Range("L19").Value = 47.4359690887577
$ws.Range("N19").Value = 37.67551094615366

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B3").Value = 168381.7769200939
$ws.Range("B4").Value = 202375.4442256468

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("C2").Value = 62730.46591140758
$ws.Range("D2").Value = 75394.77333896644
$ws.Range("E2").Value = 91976.24205358136
$ws.Range("F2").Value = 91976.24205358134
$ws.Range("G2").Value = 91976.24205358134
$ws.Range("H2").Value = 91976.24205358134
$ws.Range("L2").Value = 91976.24205358136
$ws.Range("O2").Value = 91976.24205358134
$ws.Range("C3").Value = 196825.9098199031
$ws.Range("D3").Value = 38236.46568336456
$ws.Range("E3").Value = 52530.53686621619
$ws.Range("C5").Value = 38339.65294307929
$ws.Range("D5").Value = 39312.96135688073
$ws.Range("B6").Value = -54153.64424660708
$ws.Range("C6").Value = -189659.5376291605
$ws.Range("D6").Value = -18712.55198267718
$ws.Range("E6").Value = 16649.76056833227
$ws.Range("F6").Value = 69180.29743454844
$ws.Range("G6").Value = 69180.29743454844
$ws.Range("H6").Value = 69180.29743454844
$ws.Range("I6").Value = 69180.29743454844
$ws.Range("J6").Value = 69180.29743454844
$ws.Range("K6").Value = 69180.29743454847
$ws.Range("L6").Value = 69180.29743454845
$ws.Range("M6").Value = 69180.29743454844
$ws.Range("N6").Value = 69180.29743454844
$ws.Range("O6").Value = 69180.29743454844
$ws.Range("P6").Value = 69180.29743454844

$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("C3").Value = 216.1492175724446
$ws.Range("D3").Value = 260.7963925174648

$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("C3").Value = 216.1492175724445
$ws.Range("D3").Value = 44.64717494502023
$ws.Range("E3").Value = 65.38503947111997

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G5").Value = 0.8689415781806812
$ws.Range("H5").Value = 8.899047937542903
$ws.Range("I5").Value = 33.49987019281074
$ws.Range("J5").Value = 73.75033027111266
$ws.Range("K5").Value = 110.5326272755009
$ws.Range("L5").Value = 137.125498098748
$ws.Range("M5").Value = 152.5785378897186
$ws.Range("N5").Value = 155.0474181487245
$ws.Range("O5").Value = 146.4068803306903
$ws.Range("P5").Value = 124.9548851193548
$ws.Range("Q5").Value = 93.8359148507591
$ws.Range("R5").Value = 54.58365141039226
$ws.Range("S5").Value = 19.80100621279229
$ws.Range("T5").Value = 3.803791758485934
$ws.Range("U5").Value = 0.06951532625445447
$ws.Range("G6").Value = 0.4649247321369563
$ws.Range("H6").Value = 4.490194123533237
$ws.Range("I6").Value = 16.00727696173293
$ws.Range("J6").Value = 43.92519146913236
$ws.Range("K6").Value = 75.07514852230842
$ws.Range("L6").Value = 100.9478020370177
$ws.Range("M6").Value = 117.8013235769823
$ws.Range("N6").Value = 120.9191740832867
$ws.Range("O6").Value = 110.6174208080415
$ws.Range("P6").Value = 88.78023240253862
$ws.Range("Q6").Value = 59.34723422857008
$ws.Range("R6").Value = 28.86611626373139
$ws.Range("S6").Value = 8.635772985087758
$ws.Range("T6").Value = 1.873972933481854
$ws.Range("U6").Value = 0.03058715343006293
$ws.Range("G7").Value = 0.3897772775896541
$ws.Range("H7").Value = 3.465474340751655
$ws.Range("I7").Value = 11.72166576605978
$ws.Range("J7").Value = 27.55725352558855
$ws.Range("K7").Value = 45.28503279632526
$ws.Range("L7").Value = 57.94925088819277
$ws.Range("M7").Value = 61.09935997707642
$ws.Range("N7").Value = 59.64655376060593
$ws.Range("O7").Value = 55.0932464723995
$ws.Range("P7").Value = 47.14179000957051
$ws.Range("Q7").Value = 32.63853185343913
$ws.Range("R7").Value = 17.52580377234936
$ws.Range("S7").Value = 6.792754919448789
$ws.Range("T7").Value = 1.665412004246704
$ws.Range("U7").Value = 0.02126057877761752
$ws.Range("G8").Value = 1.048427708612923
$ws.Range("H8").Value = 10.7372102708321
$ws.Range("I8").Value = 40.41950923629976
$ws.Range("J8").Value = 88.98399123388617
$ws.Range("K8").Value = 133.3639361394712
$ws.Range("L8").Value = 165.449755626934
$ws.Range("M8").Value = 184.0947318899791
$ws.Range("N8").Value = 187.0735771170756
$ws.Range("O8").Value = 176.6482740895558
$ws.Range("P8").Value = 150.7652150331742
$ws.Range("Q8").Value = 113.2183977184739
$ws.Range("R8").Value = 65.8582970511566
$ws.Range("S8").Value = 23.89104641001701
$ws.Range("T8").Value = 4.589492294453073
$ws.Range("U8").Value = 0.08387421668903385
$ws.Range("G9").Value = 0.560958278245113
$ws.Range("H9").Value = 5.417676003051488
$ws.Range("I9").Value = 19.31369510624622
$ws.Range("J9").Value = 52.99825557753501
$ws.Range("K9").Value = 90.58246021995758
$ws.Range("L9").Value = 121.799296335633
$ws.Range("M9").Value = 142.1340339220183
$ws.Range("O9").Value = 133.4662443858011
$ws.Range("P9").Value = 107.1184277115602
$ws.Range("Q9").Value = 71.60583214932356
$ws.Range("R9").Value = 34.82862011770975
$ws.Range("S9").Value = 10.41955398407041
$ws.Range("T9").Value = 2.261055516259907
$ws.Range("U9").Value = 0.03690514988454693
$ws.Range("G10").Value = 0.4702885766708382
$ws.Range("H10").Value = 4.181292981673455
$ws.Range("I10").Value = 14.14286010570121
$ws.Range("J10").Value = 33.24940237062826
$ws.Range("K10").Value = 54.63898190775737
$ws.Range("L10").Value = 69.91908529886263
$ws.Range("M10").Value = 73.71987206850237
$ws.Range("N10").Value = 71.96697828272931
$ws.Range("O10").Value = 66.47315263707451
$ws.Range("P10").Value = 56.87926567298936
$ws.Range("Q10").Value = 39.38025527013719
$ws.Range("R10").Value = 21.14588454739968
$ws.Range("S10").Value = 8.195847286163604
$ws.Range("T10").Value = 2.009414827593581
$ws.Range("U10").Value = 0.02565210418204575
$ws.Range("M11").Value = 230.2496698278013
$ws.Range("J13").Value = 41.58545896024956
$ws.Range("L13").Value = 87.4487071924806
$ws.Range("N13").Value = 90.01003351907951
$ws.Range("M17").Value = 230.2496698278014
$ws.Range("J19").Value = 41.58545896024957
$ws.Range("L19").Value = 87.44870719248058
$ws.Range("N19").Value = 90.01003351907953
